# Add the 2022-Q4 quarter sheet (and its row in 总计) to the workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" worksheet.
#    We duplicate the current "2022-Q3" sheet (2nd tab) so the new sheet
#    starts out with identical layout/styles/header row, then place it
#    right after "总计" (i.e. before the existing "2022-Q3" tab), and
#    finally overwrite its data with the 2022-Q4 numbers.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3) | Out-Null

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Drop the now-superfluous rows 5-9 that were copied from 2022-Q3 (the
# new quarter only has 3 fund rows).
$newSheet.Range("A5:H9").Delete() | Out-Null

# Helper: write a *text* value into a cell without Excel "helpfully"
# re-interpreting numeric-looking strings (e.g. "0.73", "002863") as
# numbers. We stage the text through a throw-away formula cell (whose
# computed result is a string) and Copy/Paste that computed value into
# the real destination, then wipe the scratch cell again.
$scratchCell = $newSheet.Range("Z100")
function Set-TextValue($sheet, $range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $scratch = $sheet.Range("Z100")
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy($range) | Out-Null
    $scratch.Clear() | Out-Null
}

# Row 2 - 002863 金信深圳成长灵活配置混合
Set-TextValue $newSheet $newSheet.Range("B2") "002863"
$newSheet.Range("C2").Value = "金信深圳成长灵活配置混合"
Set-TextValue $newSheet $newSheet.Range("D2") "0.73"
Set-TextValue $newSheet $newSheet.Range("E2") "91.56"
Set-TextValue $newSheet $newSheet.Range("F2") "5.15"
Set-TextValue $newSheet $newSheet.Range("G2") "0.0376"
$newSheet.Range("H2").Value = 5

# Row 3 - 014246 大摩现代服务业混合A
Set-TextValue $newSheet $newSheet.Range("B3") "014246"
$newSheet.Range("C3").Value = "大摩现代服务业混合A"
Set-TextValue $newSheet $newSheet.Range("D3") "0.17"
Set-TextValue $newSheet $newSheet.Range("E3") "86.98"
Set-TextValue $newSheet $newSheet.Range("F3") "5.87"
Set-TextValue $newSheet $newSheet.Range("G3") "0.0100"
$newSheet.Range("H3").Value = 9

# Row 4 - 014247 大摩现代服务业混合C
Set-TextValue $newSheet $newSheet.Range("B4") "014247"
$newSheet.Range("C4").Value = "大摩现代服务业混合C"
Set-TextValue $newSheet $newSheet.Range("D4") "0.06"
Set-TextValue $newSheet $newSheet.Range("E4") "86.98"
Set-TextValue $newSheet $newSheet.Range("F4") "5.87"
Set-TextValue $newSheet $newSheet.Range("G4") "0.0035"
$newSheet.Range("H4").Value = 9

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert the 2022-Q4 row at the top
#    of the data (row 2) and shift the rest of the quarters down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Duplicate row 6's formatting down into row 7 so the new last row keeps
# the same per-column styling (column A uses style index 2) as the rest
# of the table.
$summary.Range("A6:D6").Copy($summary.Range("A7:D7")) | Out-Null

# Now rewrite all the data rows (2-7) with their final values.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.05

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 8
$summary.Range("D3").Value = 0.14

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 7
$summary.Range("D4").Value = 1.22

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 2
$summary.Range("D5").Value = 0.12

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q4"
$summary.Range("C6").Value = 7
$summary.Range("D6").Value = 0.62

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q3"
$summary.Range("C7").Value = 4
$summary.Range("D7").Value = 0.08

# Restore the original active sheet (总计 / the first tab).
$summary.Activate()

Write-Output "Workbook updated."
